# Update the label in B1 from "АТА&Pos" to "Position"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "Position"

# Move the active selection to F9 (as recorded in the saved view state)
$ws.Range("F9").Select()
